# Update "Regra repasse" (col K) and "Comissão Vendedor" (col L) on the
# "Transações" sheet per the new repasse rules, then zero out the dependent
# per-client (sheet "Por Cliente", col E) and grand-total (sheet "Totais",
# cell B3) rollups that were pasted-in as static values.

$wb = $excel.ActiveWorkbook

$wsTrans = $wb.Worksheets.Item("Transações")
$wsCliente = $wb.Worksheets.Item("Por Cliente")
$wsTotais = $wb.Worksheets.Item("Totais")

# Rows whose rule was "*0,8517/2" -> becomes "*0,8517*0,25"; their
# "Comissão Vendedor" (L) result is zeroed out.
$rowsHalfToQuarter = @(
    2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,
    29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,
    53,54,55,56,57,58,59,62,63,64,65,66,67,68,69,95,96,97,98,99,100,134,135,
    136,137,138,139,140,141,142,143,144,145,146,160,161,162,163,164,165,166,
    167,168,169,170,171,172,173,174,175,177,178,179,180,181,182,183,184,185,
    186,187,188,189,190,191,192,193,194,195,196,197,198,199
)

foreach ($r in $rowsHalfToQuarter) {
    $wsTrans.Cells.Item($r, 11).Value = "*0,8517*0,25"
    $wsTrans.Cells.Item($r, 12).Value = 0
}

# Rows whose rule was "*0,8517*0,2" -> becomes "*0,8517*0,1". These rows'
# "Comissão Vendedor" (L) was already 0 and stays 0.
$rowsFifthToTenth = @(
    60,61,70,71,72,73,74,75,76,77,78,79,80,81,82,83,84,85,86,87,88,89,90,91,
    92,93,94,101,102,103,104,105,106,107,108,109,110,111,112,113,114,115,116,
    117,118,119,120,121,122,123,124,125,126,127,128,129,130,131,132,133,147,
    148,149,150,151,152,153,154,155,156,157,158,159,176
)

foreach ($r in $rowsFifthToTenth) {
    $wsTrans.Cells.Item($r, 11).Value = "*0,8517*0,1"
}

# "Por Cliente" sheet: zero out the "Comissão Vendedor" (col E) rollup for
# every client row that depended on the now-zeroed "Transações" amounts.
$rowsClienteToZero = @(2,4,5,6,8,9,14,15,16,17,18,19,20,22,23,25,26,27,28,29,30,31)

foreach ($r in $rowsClienteToZero) {
    $wsCliente.Cells.Item($r, 5).Value = 0
}

# "Totais" sheet: the grand total "Total Comissão Vendedor" also collapses
# to 0.
$wsTotais.Range("B3").Value = 0
